# Re-sort the country list (A2:B20) alphabetically by country name,
# keeping each country's "Total Projects" value intact.
# Data supplemented from funding tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Argentina", 18),
    @("Aruba", 1),
    @("Bolivia", 3),
    @("Brazil", 201),
    @("Canada", 1640),
    @("Chile", 1),
    @("Colombia", 8),
    @("Curacao", 1),
    @("Ecuador", 3),
    @("Guatemala", 1),
    @("Guyana", 1),
    @("Haiti", 3),
    @("Honduras", 1),
    @("Mexico", 9),
    @("Nicaragua", 2),
    @("Panama", 1),
    @("Peru", 11),
    @("St. Maarten", 1),
    @("USA", 2793)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
